$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A5 value
$ws.Range("A5").Value = [double]"4.1979166666666672E-2"

# Add new rows 16-20 to column A (mirrors style of A2:A15, time values)
$ws.Range("A16").Value = [double]"3.9548611111111111E-2"
$ws.Range("A17").Value = [double]"8.621527777777778E-2"
$ws.Range("A18").Value = [double]"8.1736111111111107E-2"
$ws.Range("A19").Value = [double]"2.9155092592592594E-2"
$ws.Range("A20").Value = [double]"3.636574074074074E-2"

# Match the time-format style (s="2") used by the rest of column A
$ws.Range("A16:A20").NumberFormat = $ws.Range("A15").NumberFormat

# Update formulas to cover the new range
$ws.Range("C2").Formula = "=SUM(A1:A20)"
$ws.Range("B3").Formula = "=SUM(A9:A20)"

# Update selection to D3
$ws.Range("D3").Select()
